$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every populated cell in columns B:G as plain
# text (coin names, URLs, price/volume/hour strings formatted by the
# scraping script) rather than native numbers. Force each edited cell
# to the Text number format before writing so Excel keeps the literal
# string instead of re-parsing "27.92" or "-5.08%" into a number.

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "7"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.08%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "7"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.321"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.11%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "7"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05849"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.91%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "7"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.713"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.83%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "7"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8646"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.66%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "7"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9053"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.22%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "7"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1427"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.18%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "7"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07179"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.25%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "7"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03183"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.43%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "7"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09223"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.65%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "7"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001550"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.29%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "One"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006035"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-94.12%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "7"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005870"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.55%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "7"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.500"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.17%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "7"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.224"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.07%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "7"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.202"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.99%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "7"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3167"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.92%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "7"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03457"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.19%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "7"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.10%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "7"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.527"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.08%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "7"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04157"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.11%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "7"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1378"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.20%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "7"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.005118"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "23.56%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "7"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.001226"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.16%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "7"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009994"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-17.38%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "7"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001936"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "34.02%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "7"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "7"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "7"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "7"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "7"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "7"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "7"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "7"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "7"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "7"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "7"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "7"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03861"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.88%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "7"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.69%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "7"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002199"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.02%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "7"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002948"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-48.62%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "7"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01098"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "18.61%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "7"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005240"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.00%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "7"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "7"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08978"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "54.83%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "7"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.02%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "7"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "7"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "7"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "7"
